$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift embarque-related headers, add morte-related headers ---
$ws.Range("I1").Value = "embarque.date.year"
$ws.Range("J1").Value = "age_at_embarque"
$ws.Range("K1").Value = "morte.date.year"

# New header cells L1:N1 should carry the same header style as the existing header row (e.g. K1)
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("L1").Value = "age_at_morte"
$ws.Range("M1").Value = "mission_time"
$ws.Range("N1").Value = "morte"

# --- Data rows 2-58 ---
# Columns I, K, N hold text (years-as-text / place names); prefix with an apostrophe so
# numeric-looking strings (years) are stored as text, matching the source data (embarque/morte
# year columns are text-typed, not numeric, in this dataset).
# Columns J, L, M hold numeric ages / durations.

# Row 2
$ws.Cells.Item(2, 11).Value = "'1579"
$ws.Cells.Item(2, 12).Value = 56
$ws.Cells.Item(2, 14).Value = "'?"
# Row 3
$ws.Cells.Item(3, 9).Value = "'1551"
$ws.Cells.Item(3, 10).Value = 31
$ws.Cells.Item(3, 11).Value = "'1571"
$ws.Cells.Item(3, 12).Value = 52
$ws.Cells.Item(3, 13).Value = 20
$ws.Cells.Item(3, 14).Value = "'Goa"
# Row 4
$ws.Cells.Item(4, 9).Value = "'1555"
$ws.Cells.Item(4, 10).Value = 35
$ws.Cells.Item(4, 11).Value = "'1583"
$ws.Cells.Item(4, 12).Value = 64
$ws.Cells.Item(4, 13).Value = 28
$ws.Cells.Item(4, 14).Value = "'Macau"
# Row 5
$ws.Cells.Item(5, 9).Value = "'1546"
$ws.Cells.Item(5, 10).Value = 31
$ws.Cells.Item(5, 11).Value = "'1583"
$ws.Cells.Item(5, 12).Value = 68
$ws.Cells.Item(5, 13).Value = 36
$ws.Cells.Item(5, 14).Value = "'Negapatam"
# Row 6
$ws.Cells.Item(6, 9).Value = "'1555"
$ws.Cells.Item(6, 10).Value = 26
$ws.Cells.Item(6, 11).Value = "'1572"
$ws.Cells.Item(6, 12).Value = 44
$ws.Cells.Item(6, 13).Value = 17
$ws.Cells.Item(6, 14).Value = "'Goa"
# Row 7
$ws.Cells.Item(7, 9).Value = "'1568"
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = "'1573"
$ws.Cells.Item(7, 12).Value = 46
$ws.Cells.Item(7, 13).Value = 5
$ws.Cells.Item(7, 14).Value = "'[A caminho do Japão]"
# Row 8
$ws.Cells.Item(8, 9).Value = "'1551"
$ws.Cells.Item(8, 10).Value = 21
$ws.Cells.Item(8, 11).Value = "'1582"
$ws.Cells.Item(8, 12).Value = 52
$ws.Cells.Item(8, 13).Value = 31
$ws.Cells.Item(8, 14).Value = "'Macau"
# Row 9
$ws.Cells.Item(9, 9).Value = "'1556"
$ws.Cells.Item(9, 10).Value = 24
$ws.Cells.Item(9, 11).ClearContents()
# Row 10
$ws.Cells.Item(10, 9).Value = "'1585"
$ws.Cells.Item(10, 10).Value = 42
$ws.Cells.Item(10, 11).Value = "'1598"
$ws.Cells.Item(10, 12).Value = 55
$ws.Cells.Item(10, 13).Value = 12
$ws.Cells.Item(10, 14).Value = "'[No mar, a caminho de Malaca]"
# Row 11
$ws.Cells.Item(11, 9).Value = "'1576"
$ws.Cells.Item(11, 10).Value = 36
$ws.Cells.Item(11, 11).Value = "'1599"
$ws.Cells.Item(11, 12).Value = 59
$ws.Cells.Item(11, 13).Value = 23
$ws.Cells.Item(11, 14).Value = "'Macau"
# Row 12
$ws.Cells.Item(12, 9).Value = "'1596"
$ws.Cells.Item(12, 10).Value = 49
$ws.Cells.Item(12, 11).Value = "'1613"
$ws.Cells.Item(12, 12).Value = 66
$ws.Cells.Item(12, 13).Value = 16
$ws.Cells.Item(12, 14).Value = "'Goa (Colégio Novo)"
# Row 13
$ws.Cells.Item(13, 9).Value = "'1589"
$ws.Cells.Item(13, 10).Value = 34
$ws.Cells.Item(13, 11).ClearContents()
# Row 14
$ws.Cells.Item(14, 9).Value = "'1583"
$ws.Cells.Item(14, 10).Value = 24
$ws.Cells.Item(14, 11).Value = "'1629"
$ws.Cells.Item(14, 12).Value = 71
$ws.Cells.Item(14, 13).Value = 46
$ws.Cells.Item(14, 14).Value = "'Macau"
# Row 15
$ws.Cells.Item(15, 9).Value = "'1586"
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 11).Value = "'1623"
$ws.Cells.Item(15, 12).Value = 57
$ws.Cells.Item(15, 13).Value = 36
$ws.Cells.Item(15, 14).Value = "'Hang-tcheou"
# Row 16
$ws.Cells.Item(16, 9).Value = "'1586"
$ws.Cells.Item(16, 10).Value = 16
$ws.Cells.Item(16, 11).Value = "'1632"
$ws.Cells.Item(16, 12).Value = 63
$ws.Cells.Item(16, 13).Value = 46
$ws.Cells.Item(16, 14).Value = "'Japão"
# Row 17
$ws.Cells.Item(17, 9).Value = "'1586"
$ws.Cells.Item(17, 10).Value = 19
$ws.Cells.Item(17, 11).Value = "'1607"
$ws.Cells.Item(17, 12).Value = 41
$ws.Cells.Item(17, 13).Value = 21
$ws.Cells.Item(17, 14).Value = "'Macau"
# Row 18
$ws.Cells.Item(18, 9).Value = "'1618"
$ws.Cells.Item(18, 10).Value = 49
$ws.Cells.Item(18, 11).Value = "'1633"
$ws.Cells.Item(18, 12).Value = 65
$ws.Cells.Item(18, 13).Value = 15
$ws.Cells.Item(18, 14).Value = "'Macau"
# Row 19
$ws.Cells.Item(19, 9).Value = "'1617"
$ws.Cells.Item(19, 10).Value = 47
$ws.Cells.Item(19, 11).Value = "'1635"
$ws.Cells.Item(19, 12).Value = 65
$ws.Cells.Item(19, 13).Value = 17
$ws.Cells.Item(19, 14).Value = "'Macau"
# Row 20
$ws.Cells.Item(20, 9).Value = "'1592"
$ws.Cells.Item(20, 10).Value = 23
$ws.Cells.Item(20, 11).Value = "'1626"
$ws.Cells.Item(20, 12).Value = 57
$ws.Cells.Item(20, 13).Value = 34
$ws.Cells.Item(20, 14).Value = "'Nagasaki"
# Row 21
$ws.Cells.Item(21, 9).Value = "'1593"
$ws.Cells.Item(21, 10).Value = 21
$ws.Cells.Item(21, 11).Value = "'1649"
$ws.Cells.Item(21, 12).Value = 78
$ws.Cells.Item(21, 13).Value = 56
$ws.Cells.Item(21, 14).Value = "'Cantão"
# Row 22
$ws.Cells.Item(22, 9).Value = "'1597"
$ws.Cells.Item(22, 10).Value = 23
$ws.Cells.Item(22, 11).Value = "'1607"
$ws.Cells.Item(22, 12).Value = 34
$ws.Cells.Item(22, 13).Value = 10
$ws.Cells.Item(22, 14).Value = "'Japão"
# Row 23
$ws.Cells.Item(23, 9).Value = "'1601"
$ws.Cells.Item(23, 10).Value = 26
$ws.Cells.Item(23, 11).Value = "'1659"
$ws.Cells.Item(23, 12).Value = 84
$ws.Cells.Item(23, 13).Value = 57
$ws.Cells.Item(23, 14).Value = "'Hangchow"
# Row 24
$ws.Cells.Item(24, 9).Value = "'1601"
$ws.Cells.Item(24, 10).Value = 21
$ws.Cells.Item(24, 11).Value = "'1614"
$ws.Cells.Item(24, 12).Value = 34
$ws.Cells.Item(24, 13).Value = 13
$ws.Cells.Item(24, 14).Value = "'Nanquim"
# Row 25
$ws.Cells.Item(25, 9).Value = "'1600"
$ws.Cells.Item(25, 10).Value = 19
$ws.Cells.Item(25, 11).Value = "'1634"
$ws.Cells.Item(25, 12).Value = 53
$ws.Cells.Item(25, 13).Value = 33
$ws.Cells.Item(25, 14).Value = "'Goa"
# Row 26
$ws.Cells.Item(26, 9).Value = "'1618"
$ws.Cells.Item(26, 10).Value = 28
$ws.Cells.Item(26, 11).Value = "'1660"
$ws.Cells.Item(26, 12).Value = 71
$ws.Cells.Item(26, 13).Value = 42
$ws.Cells.Item(26, 14).Value = "'Macau"
# Row 27
$ws.Cells.Item(27, 9).Value = "'1623"
$ws.Cells.Item(27, 10).Value = 30
$ws.Cells.Item(27, 11).Value = "'1677"
$ws.Cells.Item(27, 12).Value = 84
$ws.Cells.Item(27, 13).Value = 53
$ws.Cells.Item(27, 14).Value = "'Foochow, Fukien"
# Row 28
$ws.Cells.Item(28, 9).Value = "'1623"
$ws.Cells.Item(28, 10).Value = 30
$ws.Cells.Item(28, 11).Value = "'1646"
$ws.Cells.Item(28, 12).Value = 53
$ws.Cells.Item(28, 13).Value = 22
$ws.Cells.Item(28, 14).Value = "'[Perto de Hainan]"
# Row 29
$ws.Cells.Item(29, 9).Value = "'1618"
$ws.Cells.Item(29, 10).Value = 28
$ws.Cells.Item(29, 11).Value = "'1653"
$ws.Cells.Item(29, 12).Value = 64
$ws.Cells.Item(29, 13).Value = 35
$ws.Cells.Item(29, 14).Value = "'Macau"
# Row 30
$ws.Cells.Item(30, 9).Value = "'1637"
$ws.Cells.Item(30, 10).Value = 39
$ws.Cells.Item(30, 11).Value = "'1664"
$ws.Cells.Item(30, 12).Value = 66
$ws.Cells.Item(30, 13).Value = 27
$ws.Cells.Item(30, 14).Value = "'Macau"
# Row 31
$ws.Cells.Item(31, 9).Value = "'1634"
$ws.Cells.Item(31, 10).Value = 23
$ws.Cells.Item(31, 11).Value = "'1677"
$ws.Cells.Item(31, 12).Value = 66
$ws.Cells.Item(31, 13).Value = 43
$ws.Cells.Item(31, 14).Value = "'Pequim"
# Row 32
$ws.Cells.Item(32, 9).Value = "'1640"
$ws.Cells.Item(32, 10).Value = 23
$ws.Cells.Item(32, 11).Value = "'1667"
$ws.Cells.Item(32, 12).Value = 50
$ws.Cells.Item(32, 13).Value = 26
$ws.Cells.Item(32, 14).Value = "'Costa da Cochinchina"
# Row 33
$ws.Cells.Item(33, 9).Value = "'1643"
$ws.Cells.Item(33, 10).Value = 25
$ws.Cells.Item(33, 11).ClearContents()
# Row 34
$ws.Cells.Item(34, 9).Value = "'1643"
$ws.Cells.Item(34, 10).Value = 21
$ws.Cells.Item(34, 11).Value = "'1677"
$ws.Cells.Item(34, 12).Value = 56
$ws.Cells.Item(34, 13).Value = 34
$ws.Cells.Item(34, 14).Value = "'Nanquim"
# Row 35
$ws.Cells.Item(35, 9).Value = "'1657"
$ws.Cells.Item(35, 10).Value = 31
$ws.Cells.Item(35, 11).Value = "'1661"
$ws.Cells.Item(35, 12).Value = 36
$ws.Cells.Item(35, 13).Value = 4
$ws.Cells.Item(35, 14).Value = "'Foochow (Fou-tcheou fou, Fukien)"
# Row 36
$ws.Cells.Item(36, 9).Value = "'1657"
$ws.Cells.Item(36, 10).Value = 26
$ws.Cells.Item(36, 11).Value = "'1657"
$ws.Cells.Item(36, 12).Value = 27
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = "'[No mar, depois do Cabo da Boa Esperança]"
# Row 37
$ws.Cells.Item(37, 9).Value = "'1666"
$ws.Cells.Item(37, 10).Value = 20
$ws.Cells.Item(37, 11).Value = "'1708"
$ws.Cells.Item(37, 12).Value = 63
$ws.Cells.Item(37, 13).Value = 42
$ws.Cells.Item(37, 14).Value = "'Pequim"
# Row 38
$ws.Cells.Item(38, 9).Value = "'1694"
$ws.Cells.Item(38, 10).Value = 44
$ws.Cells.Item(38, 11).Value = "'1709"
$ws.Cells.Item(38, 12).Value = 58
$ws.Cells.Item(38, 13).Value = 14
$ws.Cells.Item(38, 14).Value = "'Goa"
# Row 39
$ws.Cells.Item(39, 11).ClearContents()
# Row 40
$ws.Cells.Item(40, 9).Value = "'1694"
$ws.Cells.Item(40, 10).Value = 40
$ws.Cells.Item(40, 11).Value = "'1726"
$ws.Cells.Item(40, 12).Value = 72
$ws.Cells.Item(40, 13).Value = 31
$ws.Cells.Item(40, 14).Value = "'Nanquim"
# Row 41
$ws.Cells.Item(41, 9).Value = "'1680"
$ws.Cells.Item(41, 10).Value = 24
$ws.Cells.Item(41, 11).Value = "'1741"
$ws.Cells.Item(41, 12).Value = 84
$ws.Cells.Item(41, 13).Value = 60
$ws.Cells.Item(41, 14).Value = "'Macau"
# Row 42
$ws.Cells.Item(42, 9).Value = "'1680"
$ws.Cells.Item(42, 10).Value = 21
$ws.Cells.Item(42, 11).Value = "'1728"
$ws.Cells.Item(42, 12).Value = 69
$ws.Cells.Item(42, 13).Value = 48
$ws.Cells.Item(42, 14).Value = "'Lisboa"
# Row 43
$ws.Cells.Item(43, 9).Value = "'1681"
$ws.Cells.Item(43, 10).Value = 18
$ws.Cells.Item(43, 11).Value = "'1731"
$ws.Cells.Item(43, 12).Value = 68
$ws.Cells.Item(43, 13).Value = 49
$ws.Cells.Item(43, 14).Value = "'Macau"
# Row 44
$ws.Cells.Item(44, 9).Value = "'1682"
$ws.Cells.Item(44, 10).Value = 24
$ws.Cells.Item(44, 11).Value = "'1730"
$ws.Cells.Item(44, 12).Value = 72
$ws.Cells.Item(44, 13).Value = 48
$ws.Cells.Item(44, 14).Value = "'Coimbra"
# Row 45
$ws.Cells.Item(45, 9).Value = "'1685"
$ws.Cells.Item(45, 10).Value = 22
$ws.Cells.Item(45, 11).Value = "'1710"
$ws.Cells.Item(45, 12).Value = 47
$ws.Cells.Item(45, 13).Value = 25
$ws.Cells.Item(45, 14).Value = "'Macau"
# Row 46
$ws.Cells.Item(46, 9).Value = "'1687"
$ws.Cells.Item(46, 10).Value = 20
$ws.Cells.Item(46, 11).Value = "'1729"
$ws.Cells.Item(46, 12).Value = 62
$ws.Cells.Item(46, 13).Value = 42
$ws.Cells.Item(46, 14).Value = "'Tonquim"
# Row 47
$ws.Cells.Item(47, 9).Value = "'1695"
$ws.Cells.Item(47, 10).Value = 26
$ws.Cells.Item(47, 11).Value = "'1699"
$ws.Cells.Item(47, 12).Value = 29
$ws.Cells.Item(47, 13).Value = 3
$ws.Cells.Item(47, 14).Value = "'Macau"
# Row 48
$ws.Cells.Item(48, 9).Value = "'1690"
$ws.Cells.Item(48, 10).Value = 19
$ws.Cells.Item(48, 11).ClearContents()
# Row 49
$ws.Cells.Item(49, 9).Value = "'1695"
$ws.Cells.Item(49, 10).Value = 19
$ws.Cells.Item(49, 11).ClearContents()
# Row 50
$ws.Cells.Item(50, 9).Value = "'1696"
$ws.Cells.Item(50, 10).Value = 19
$ws.Cells.Item(50, 11).Value = "'1735"
$ws.Cells.Item(50, 12).Value = 57
$ws.Cells.Item(50, 13).Value = 38
$ws.Cells.Item(50, 14).Value = "'Pequim"
# Row 51
$ws.Cells.Item(51, 9).Value = "'1694"
$ws.Cells.Item(51, 10).Value = 19
$ws.Cells.Item(51, 11).Value = "'1721"
$ws.Cells.Item(51, 12).Value = 46
$ws.Cells.Item(51, 13).Value = 26
$ws.Cells.Item(51, 14).Value = "'Goa"
# Row 52
$ws.Cells.Item(52, 9).Value = "'1695"
$ws.Cells.Item(52, 10).Value = 18
$ws.Cells.Item(52, 11).Value = "'1734"
$ws.Cells.Item(52, 12).Value = 57
$ws.Cells.Item(52, 13).Value = 38
$ws.Cells.Item(52, 14).Value = "'?"
# Row 53
$ws.Cells.Item(53, 9).Value = "'1695"
$ws.Cells.Item(53, 10).Value = 20
$ws.Cells.Item(53, 11).ClearContents()
# Row 54
$ws.Cells.Item(54, 9).Value = "'1715"
$ws.Cells.Item(54, 10).Value = 33
$ws.Cells.Item(54, 11).Value = "'1752"
$ws.Cells.Item(54, 12).Value = 70
$ws.Cells.Item(54, 13).Value = 37
$ws.Cells.Item(54, 14).Value = "'?"
# Row 55
$ws.Cells.Item(55, 9).Value = "'1714"
$ws.Cells.Item(55, 10).Value = 24
$ws.Cells.Item(55, 11).ClearContents()
# Row 56
$ws.Cells.Item(56, 9).Value = "'1742"
$ws.Cells.Item(56, 10).Value = 34
$ws.Cells.Item(56, 11).Value = "'1764"
$ws.Cells.Item(56, 12).Value = 56
$ws.Cells.Item(56, 13).Value = 22
$ws.Cells.Item(56, 14).Value = "'Castel Gandolfo"
# Row 57
$ws.Cells.Item(57, 9).Value = "'1727"
$ws.Cells.Item(57, 10).Value = 20
$ws.Cells.Item(57, 11).Value = "'1751"
$ws.Cells.Item(57, 12).Value = 44
$ws.Cells.Item(57, 13).Value = 23
$ws.Cells.Item(57, 14).Value = "'Pequim"
# Row 58
$ws.Cells.Item(58, 9).Value = "'1750"
$ws.Cells.Item(58, 10).Value = 25
$ws.Cells.Item(58, 11).Value = "'1776"
$ws.Cells.Item(58, 12).Value = 51
$ws.Cells.Item(58, 13).Value = 26
$ws.Cells.Item(58, 14).Value = "'?"

Write-Output "applied edit"
